# Validation case now includes graphs for K and P
#
# The "All same / Very different / Somewhat different" scenario headers
# become single, underscore-joined tokens (so they can be used as
# chart/series names for the new K & P graphs), and the header cells for
# those three columns drop their bordered header style. The active
# selection is moved onto the new header range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the scenario headers: spaces -> underscores
$ws.Range("I1").Value = "All_same"
$ws.Range("J1").Value = "Very_different"
$ws.Range("K1").Value = "Somewhat_different"

# Those header cells revert to the default (un-bordered) style
$ws.Range("I1:K1").Style = "Normal"

# Select the new header range (I1 active, I1:K1 highlighted)
$ws.Range("I1:K1").Select() | Out-Null
